# Adds a new "Duplicated Columns" worksheet (a copy of "Data OK" with an
# extra duplicated "Int Column" appended as column L), and updates the
# active sheet/selection bookkeeping accordingly.

$wb = $excel.ActiveWorkbook

# "Data OK" is the first worksheet and is the source of the data we
# duplicate into the new sheet.
$sourceSheet = $wb.Worksheets.Item(1)

# Selection on "Data OK" moves from K3 to K2, and it stops being the
# active/selected tab (the new sheet will become active instead).
$sourceSheet.Range("K2").Select()

# Add the new worksheet as the last tab in the workbook.
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "Duplicated Columns"

# Copy over the cell formatting (number formats / fonts) from the source
# sheet so the same style indexes (header style, date style, etc.) get
# reused instead of new styles being created.
$sourceSheet.Range("A1:K6").Copy()
$newSheet.Range("A1").PasteSpecial(-4122) # xlPasteFormats

# Column L duplicates column A ("Int Column") - copy its formatting too.
$sourceSheet.Range("A1:A6").Copy()
$newSheet.Range("L1").PasteSpecial(-4122) # xlPasteFormats

# Headers (row 1)
$newSheet.Range("A1").Value = "Int Column"
$newSheet.Range("B1").Value = "Decimal Column"
$newSheet.Range("C1").Value = "Float Column"
$newSheet.Range("D1").Value = "Nullable Int Column"
$newSheet.Range("E1").Value = "Nullable Decimal Column"
$newSheet.Range("F1").Value = "Nullable Float Column"
$newSheet.Range("G1").Value = "String Column"
$newSheet.Range("H1").Value = "Date Column"
$newSheet.Range("I1").Value = "Nullable Date Column"
$newSheet.Range("J1").Value = "Boolean column"
$newSheet.Range("K1").Value = "Nullable Boolean column"
$newSheet.Range("L1").Value = "Int Column"

# Row 2
$newSheet.Range("A2").Value = 1
$newSheet.Range("B2").Value = 1.25
$newSheet.Range("C2").Value = 1.25
$newSheet.Range("D2").Value = 1
$newSheet.Range("E2").Value = 1.25
$newSheet.Range("F2").Value = 1.25
$newSheet.Range("G2").Value = "Item 1"
$newSheet.Range("H2").Value = 36526
$newSheet.Range("I2").Value = 36526
$newSheet.Range("J2").Value = 1
$newSheet.Range("K2").Value = 1
$newSheet.Range("L2").Value = 33

# Row 3
$newSheet.Range("A3").Value = 2
$newSheet.Range("B3").Value = 2.25
$newSheet.Range("C3").Value = 2.25
$newSheet.Range("G3").Value = "Item 2"
$newSheet.Range("H3").Value = 36527
$newSheet.Range("J3").Value = "Y"
$newSheet.Range("L3").Value = 33

# Row 4
$newSheet.Range("A4").Value = 3
$newSheet.Range("B4").Value = 3.75
$newSheet.Range("C4").Value = 3.75
$newSheet.Range("D4").Value = 3
$newSheet.Range("E4").Value = 3.75
$newSheet.Range("F4").Value = 3.75
$newSheet.Range("G4").Value = "Item 3"
$newSheet.Range("H4").Value = 36528
$newSheet.Range("I4").Value = 36528
$newSheet.Range("J4").Value = 0
$newSheet.Range("K4").Value = 0
$newSheet.Range("L4").Value = 33

# Row 5
$newSheet.Range("A5").Value = 4
$newSheet.Range("B5").Value = 4.25
$newSheet.Range("C5").Value = 4.25
$newSheet.Range("G5").Value = "Item 4"
$newSheet.Range("H5").Value = 36529
$newSheet.Range("J5").Value = "N"
$newSheet.Range("L5").Value = 33

# Row 6
$newSheet.Range("A6").Value = 5
$newSheet.Range("B6").Value = 5
$newSheet.Range("C6").Value = 5
$newSheet.Range("D6").Value = 6
$newSheet.Range("E6").Value = 5
$newSheet.Range("F6").Value = 5
$newSheet.Range("G6").Value = "Item 5"
$newSheet.Range("H6").Value = 36530
$newSheet.Range("I6").Value = 36530
$newSheet.Range("J6").Value = "S"
$newSheet.Range("K6").Value = "S"
$newSheet.Range("L6").Value = 33

# New sheet becomes the active/selected tab, with B2 selected.
$newSheet.Range("B2").Select()
